$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scheduled-task refresh: re-stamp A16's timestamp with fuller precision
$ws.Cells.Item(16, 1).Value = 45866.70861427084

# Capture the date/time number format used by existing readings in column A
$dateFormat = $ws.Cells.Item(16, 1).NumberFormat

# Append the newest sensor reading as row 17
$ws.Cells.Item(17, 1).Value = 45866.75030248108
$ws.Cells.Item(17, 1).NumberFormat = $dateFormat
$ws.Cells.Item(17, 2).Value = 2025
$ws.Cells.Item(17, 3).Value = 31
$ws.Cells.Item(17, 4).Value = 17.83
$ws.Cells.Item(17, 5).Value = 77.62
$ws.Cells.Item(17, 6).Value = 8.789999999999999
$ws.Cells.Item(17, 7).Value = 5.88
$ws.Cells.Item(17, 8).Value = "ESE"
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = "18:00:26"
